$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status roll-up: every cell still showing the placeholder "Ready for
#    handoff" status is now reporting the finished handback state.
# ---------------------------------------------------------------------------
$handedBack = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $handedBack
$wsOverview.Range("C2").Value = $handedBack
$wsOverview.Range("B3").Value = $handedBack
$wsOverview.Range("C3").Value = $handedBack

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $handedBack
$wsZhCn.Range("C3").Value = $handedBack

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $handedBack
$wsDeDe.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" (F) / "Latest Handback File"
#    (G) with hyperlinked file names, and stamp the real handback datetime
#    into "Latest Handback DateTime" (H), replacing the 0001-01-01 sentinel.
# ---------------------------------------------------------------------------
$zhMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/9159e8dfa2a32c6750fe276bc15b34e079760b04/e2e/56e93842-d61c-43bc-ae2e-b9f0d9872459.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/55232264a51ee30d2d61e1d871ba5b0bcb57320a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/56e93842-d61c-43bc-ae2e-b9f0d9872459.fb070a450de1b107743456862a49a818a996f52c.zh-cn.xlf"
$zhMdName  = "56e93842-d61c-43bc-ae2e-b9f0d9872459.md"
$zhXlfName = "56e93842-d61c-43bc-ae2e-b9f0d9872459.fb070a450de1b107743456862a49a818a996f52c.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhMdUrl, "", "", $zhMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhXlfUrl, "", "", $zhXlfName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhMdUrl, "", "", $zhMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhXlfUrl, "", "", $zhXlfName)

$wsZhCn.Range("H2").Value = "2016-03-18 04:34:54"
$wsZhCn.Range("H3").Value = "2016-03-18 04:34:54"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same treatment, with the de-de target/handback file names
#    and its own (later) handback timestamp.
# ---------------------------------------------------------------------------
$deMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/9159e8dfa2a32c6750fe276bc15b34e079760b04/e2e/56e93842-d61c-43bc-ae2e-b9f0d9872459.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8317ed5a3f040ef10d1fccc311e28de1d98612e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/56e93842-d61c-43bc-ae2e-b9f0d9872459.fb070a450de1b107743456862a49a818a996f52c.de-de.xlf"
$deMdName  = "56e93842-d61c-43bc-ae2e-b9f0d9872459.md"
$deXlfName = "56e93842-d61c-43bc-ae2e-b9f0d9872459.fb070a450de1b107743456862a49a818a996f52c.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deMdUrl, "", "", $deMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deXlfUrl, "", "", $deXlfName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deMdUrl, "", "", $deMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deXlfUrl, "", "", $deXlfName)

$wsDeDe.Range("H2").Value = "2016-03-18 04:34:59"
$wsDeDe.Range("H3").Value = "2016-03-18 04:34:59"
